$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new date columns before the old column C (shifts old C and
# everything right of it over to E). The new C/D cells inherit column B's
# formatting from the insert itself.
$ws.Columns("C:D").Insert()

# Preserve the header value that was riding in B1 ("Jun_13") by relocating
# it to D1, then stamp in the two new period headers.
$ws.Range("D1").Value = $ws.Range("B1").Value2
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Preserve the highlighted rating that was riding in B5 by relocating it to
# D5 (D5 already carries the inherited highlight fill from the column
# insert, so only the value needs to move).
$ws.Range("D5").Value = $ws.Range("B5").Value2

# Fill the two new columns with "UN" for every data row, matching column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    if ($r -ne 5) {
        $ws.Cells.Item($r, 4).Value = "UN"
    }
}

# Row 5's B and C cells inherited the yellow highlight from the insert;
# reset them back to the workbook's default (unhighlighted) style.
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Style = "Normal"
$ws.Range("B5").Value = "UN"
